$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the specific Price (D) / Volume (E) cells being refreshed to Text format
# so values like "1.00", "165.00", or "3.438.74" are stored verbatim instead of
# being auto-parsed as numbers/percentages by Excel.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E51").NumberFormat = "@"

$ws.Range("D2").Value = "64.281.92"
$ws.Range("E2").Value = "  -1.05%  "
$ws.Range("D3").Value = "3.438.74"
$ws.Range("E3").Value = "  -0.02%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "573.68"
$ws.Range("E5").Value = "  -0.28%  "
$ws.Range("D6").Value = "165.00"
$ws.Range("E6").Value = "  +3.41%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "3.440.63"
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").Value = "0.555"
$ws.Range("E9").Value = "  -5.78%  "
$ws.Range("E10").Value = "  +0.73%  "
$ws.Range("D11").Value = "0.120"
$ws.Range("E11").Value = "  -1.88%  "
$ws.Range("E12").Value = "  -4.72%  "
$ws.Range("D13").Value = "4.035.03"
$ws.Range("E13").Value = "  +0.08%  "
$ws.Range("E14").Value = "  +1.31%  "
$ws.Range("D15").Value = "27.40"
$ws.Range("E16").Value = "  -7.20%  "
$ws.Range("D17").Value = "64.351.10"
$ws.Range("E17").Value = "  -0.98%  "
$ws.Range("D18").Value = "3.418.94"
$ws.Range("E18").Value = "  -0.46%  "
$ws.Range("D19").Value = "6.16"
$ws.Range("E19").Value = "  -3.42%  "
$ws.Range("D20").Value = "13.72"
$ws.Range("E20").Value = "  -1.56%  "
$ws.Range("D21").Value = "379.87"
$ws.Range("E21").Value = "  -0.85%  "
$ws.Range("E22").Value = "  -1.92%  "
$ws.Range("E23").Value = "  -0.16%  "
$ws.Range("D24").Value = "71.59"
$ws.Range("E24").Value = "  -0.77%  "
$ws.Range("D25").Value = "0.522"
$ws.Range("E25").Value = "  -5.22%  "
$ws.Range("E26").Value = "  -1.91%  "
$ws.Range("D27").Value = "9.62"
$ws.Range("E27").Value = "  -3.14%  "
$ws.Range("E28").Value = "  +0.19%  "
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.21%  "
$ws.Range("D30").Value = "6.11"
$ws.Range("E30").Value = "  -0.63%  "
$ws.Range("D31").Value = "1.41"
$ws.Range("E31").Value = "  -4.76%  "
$ws.Range("E32").Value = "  +0.10%  "
$ws.Range("D33").Value = "23.06"
$ws.Range("E33").Value = "  -1.13%  "
$ws.Range("E34").Value = "  +1.05%  "
$ws.Range("E35").Value = "  -3.98%  "
$ws.Range("D36").Value = "160.04"
$ws.Range("E36").Value = "  -0.56%  "
$ws.Range("D37").Value = "0.864"
$ws.Range("E37").Value = "  +11.37%  "
$ws.Range("E38").Value = "  -4.23%  "
$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").Value = "2.824.23"
$ws.Range("E39").Value = "  -2.98%  "
$ws.Range("B40").Value = "EnergySwap"
$ws.Range("C40").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D40").Value = "26.24"
$ws.Range("E40").Value = "  -0.83%  "
$ws.Range("E41").Value = "  -2.74%  "
$ws.Range("D42").Value = "26.60"
$ws.Range("E42").Value = "  +1.72%  "
$ws.Range("D43").Value = "43.01"
$ws.Range("E43").Value = "  -0.75%  "
$ws.Range("D44").Value = "6.48"
$ws.Range("E44").Value = "  -4.80%  "
$ws.Range("E45").Value = "  -3.30%  "
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").Value = "0.0310"
$ws.Range("E46").Value = "  -2.55%  "
$ws.Range("B47").Value = "dogwifhat"
$ws.Range("C47").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D47").Value = "2.50"
$ws.Range("E47").Value = "  +9.73%  "
$ws.Range("D48").Value = "335.04"
$ws.Range("E48").Value = "  +5.90%  "
$ws.Range("E49").Value = "  -2.26%  "
$ws.Range("E50").Value = "  -2.84%  "
$ws.Range("E51").Value = "  -3.07%  "
